$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the "Lasso Regression+normalization+ lag1" model result
# right after the existing "Lasso Regression+normalization" row (row 5), pushing
# the remaining rows down by one.
$ws.Rows("6:6").Insert()

# Copy the formatting from the row that just got pushed down (now row 7, which
# used to be row 6) into the newly inserted blank row so it keeps the same
# bordered-table look as the rest of the data rows.
$ws.Range("A7:C7").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122)

# Fill in the new row's data.
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "Lasso Regression+normalization+ lag1"
$ws.Cells.Item(6, 3).Value = 81.344944740947696

# Renumber the "Id" column for all the rows that shifted down so the sequence
# stays contiguous (1..9 instead of skipping from 4 to 5-8).
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(11, 1).Value = 9

# Match the author's final selection/active cell.
$ws.Range("E14").Select()
